$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; Col = 2; Value = 1.02 },
    @{ Row = 2; Col = 3; Value = 1.019176727591457 },
    @{ Row = 2; Col = 4; Value = 1.025771812156901 },
    @{ Row = 2; Col = 5; Value = 1.047343321621769 },
    @{ Row = 2; Col = 6; Value = 1.0512274978723 },
    @{ Row = 2; Col = 9; Value = 1.030782939047553 },
    @{ Row = 2; Col = 10; Value = 1.024380974140494 },
    @{ Row = 2; Col = 11; Value = 1.028596406691454 },
    @{ Row = 2; Col = 12; Value = 1.050106224282141 },
    @{ Row = 2; Col = 13; Value = 1.053979583971713 },
    @{ Row = 2; Col = 14; Value = 1.012228640839454 },
    @{ Row = 3; Col = 2; Value = 1.02 },
    @{ Row = 3; Col = 3; Value = 1.02001514217408 },
    @{ Row = 3; Col = 4; Value = 1.026395135490157 },
    @{ Row = 3; Col = 5; Value = 1.048396122432111 },
    @{ Row = 3; Col = 6; Value = 1.052298800540129 },
    @{ Row = 3; Col = 9; Value = 1.03093890399232 },
    @{ Row = 3; Col = 10; Value = 1.024856611728966 },
    @{ Row = 3; Col = 11; Value = 1.029027603041706 },
    @{ Row = 3; Col = 12; Value = 1.050970225228418 },
    @{ Row = 3; Col = 13; Value = 1.05486282527713 },
    @{ Row = 3; Col = 14; Value = 1.012385603675983 },
    @{ Row = 4; Col = 2; Value = 1.02 },
    @{ Row = 4; Col = 3; Value = 1.020557777247906 },
    @{ Row = 4; Col = 4; Value = 1.026798146042168 },
    @{ Row = 4; Col = 5; Value = 1.049078273633288 },
    @{ Row = 4; Col = 6; Value = 1.052992612134697 },
    @{ Row = 4; Col = 9; Value = 1.031037988370779 },
    @{ Row = 4; Col = 10; Value = 1.025163875167154 },
    @{ Row = 4; Col = 11; Value = 1.029305606204582 },
    @{ Row = 4; Col = 12; Value = 1.051529608494148 },
    @{ Row = 4; Col = 13; Value = 1.055434350388573 },
    @{ Row = 4; Col = 14; Value = 1.012486989920407 },
    @{ Row = 5; Col = 2; Value = 1.02 },
    @{ Row = 5; Col = 3; Value = 1.020785929681738 },
    @{ Row = 5; Col = 4; Value = 1.026967493171878 },
    @{ Row = 5; Col = 5; Value = 1.049365269218616 },
    @{ Row = 5; Col = 6; Value = 1.053284435072504 },
    @{ Row = 5; Col = 9; Value = 1.031079203290383 },
    @{ Row = 5; Col = 10; Value = 1.025292926614376 },
    @{ Row = 5; Col = 11; Value = 1.029422235845268 },
    @{ Row = 5; Col = 12; Value = 1.05176484841891 },
    @{ Row = 5; Col = 13; Value = 1.05567462088453 },
    @{ Row = 5; Col = 14; Value = 1.012529569444514 },
    @{ Row = 6; Col = 2; Value = 1.02 },
    @{ Row = 6; Col = 3; Value = 1.020824239123682 },
    @{ Row = 6; Col = 4; Value = 1.026995922650346 },
    @{ Row = 6; Col = 5; Value = 1.049413469874743 },
    @{ Row = 6; Col = 6; Value = 1.053333441880559 },
    @{ Row = 6; Col = 9; Value = 1.031086097623704 },
    @{ Row = 6; Col = 10; Value = 1.025314587731162 },
    @{ Row = 6; Col = 11; Value = 1.029441804198072 },
    @{ Row = 6; Col = 12; Value = 1.051804350627935 },
    @{ Row = 6; Col = 13; Value = 1.055714963413382 },
    @{ Row = 6; Col = 14; Value = 1.012536716187365 },
    @{ Row = 7; Col = 2; Value = 1.02 },
    @{ Row = 7; Col = 3; Value = 1.02056082572086 },
    @{ Row = 7; Col = 4; Value = 1.026800409177283 },
    @{ Row = 7; Col = 5; Value = 1.049082107621852 },
    @{ Row = 7; Col = 6; Value = 1.05299651092013 },
    @{ Row = 7; Col = 9; Value = 1.031038540817018 },
    @{ Row = 7; Col = 10; Value = 1.025165600039939 },
    @{ Row = 7; Col = 11; Value = 1.02930716557141 },
    @{ Row = 7; Col = 12; Value = 1.051532751486883 },
    @{ Row = 7; Col = 13; Value = 1.055437560890443 },
    @{ Row = 7; Col = 14; Value = 1.01248755904049 },
    @{ Row = 8; Col = 2; Value = 1.02 },
    @{ Row = 8; Col = 3; Value = 1.019460047155285 },
    @{ Row = 8; Col = 4; Value = 1.025982532780721 },
    @{ Row = 8; Col = 5; Value = 1.047698929937391 },
    @{ Row = 8; Col = 6; Value = 1.051589423429744 },
    @{ Row = 8; Col = 9; Value = 1.030836027699847 },
    @{ Row = 8; Col = 10; Value = 1.024541822242671 },
    @{ Row = 8; Col = 11; Value = 1.028742339886961 },
    @{ Row = 8; Col = 12; Value = 1.050398151437204 },
    @{ Row = 8; Col = 13; Value = 1.054278077178519 },
    @{ Row = 8; Col = 14; Value = 1.012281724060586 },
    @{ Row = 9; Col = 2; Value = 1.02 },
    @{ Row = 9; Col = 3; Value = 1.017521344025282 },
    @{ Row = 9; Col = 4; Value = 1.024538938480866 },
    @{ Row = 9; Col = 5; Value = 1.045268661361732 },
    @{ Row = 9; Col = 6; Value = 1.049114640656994 },
    @{ Row = 9; Col = 9; Value = 1.030465150318253 },
    @{ Row = 9; Col = 10; Value = 1.02343881509686 },
    @{ Row = 9; Col = 11; Value = 1.027739362661709 },
    @{ Row = 9; Col = 12; Value = 1.048401297093579 },
    @{ Row = 9; Col = 13; Value = 1.052235022035715 },
    @{ Row = 9; Col = 14; Value = 1.011917659405444 },
    @{ Row = 10; Col = 2; Value = 1.02 },
    @{ Row = 10; Col = 3; Value = 1.016229624949635 },
    @{ Row = 10; Col = 4; Value = 1.023575026922713 },
    @{ Row = 10; Col = 5; Value = 1.043653274077265 },
    @{ Row = 10; Col = 6; Value = 1.0474679809078 },
    @{ Row = 10; Col = 9; Value = 1.03020851442058 },
    @{ Row = 10; Col = 10; Value = 1.022700959051379 },
    @{ Row = 10; Col = 11; Value = 1.027065613868519 },
    @{ Row = 10; Col = 12; Value = 1.04707174942545 },
    @{ Row = 10; Col = 13; Value = 1.05087309997167 },
    @{ Row = 10; Col = 14; Value = 1.011674057644949 },
    @{ Row = 11; Col = 2; Value = 1.02 },
    @{ Row = 11; Col = 3; Value = 1.015670488692975 },
    @{ Row = 11; Col = 4; Value = 1.023157302406753 },
    @{ Row = 11; Col = 5; Value = 1.042954937925285 },
    @{ Row = 11; Col = 6; Value = 1.046755725854209 },
    @{ Row = 11; Col = 9; Value = 1.030095171494241 },
    @{ Row = 11; Col = 10; Value = 1.022380873916253 },
    @{ Row = 11; Col = 11; Value = 1.026772678184368 },
    @{ Row = 11; Col = 12; Value = 1.046496448422882 },
    @{ Row = 11; Col = 13; Value = 1.050283408825275 },
    @{ Row = 11; Col = 14; Value = 1.011568368050634 },
    @{ Row = 12; Col = 2; Value = 1.02 },
    @{ Row = 12; Col = 3; Value = 1.015462829752258 },
    @{ Row = 12; Col = 4; Value = 1.023002090693726 },
    @{ Row = 12; Col = 5; Value = 1.042695716197299 },
    @{ Row = 12; Col = 6; Value = 1.046491277354493 },
    @{ Row = 12; Col = 9; Value = 1.030052738359189 },
    @{ Row = 12; Col = 10; Value = 1.022261892537724 },
    @{ Row = 12; Col = 11; Value = 1.026663689927401 },
    @{ Row = 12; Col = 12; Value = 1.046282816919661 },
    @{ Row = 12; Col = 13; Value = 1.050064376444838 },
    @{ Row = 12; Col = 14; Value = 1.011529079244975 },
    @{ Row = 13; Col = 2; Value = 1.02 },
    @{ Row = 13; Col = 3; Value = 1.015507371962393 },
    @{ Row = 13; Col = 4; Value = 1.023035386381669 },
    @{ Row = 13; Col = 5; Value = 1.042751312372114 },
    @{ Row = 13; Col = 6; Value = 1.046547997252901 },
    @{ Row = 13; Col = 9; Value = 1.030061855447199 },
    @{ Row = 13; Col = 10; Value = 1.022287418392443 },
    @{ Row = 13; Col = 11; Value = 1.026687076356755 },
    @{ Row = 13; Col = 12; Value = 1.046328638837846 },
    @{ Row = 13; Col = 13; Value = 1.050111359389757 },
    @{ Row = 13; Col = 14; Value = 1.011537508224007 },
    @{ Row = 14; Col = 2; Value = 1.02 },
    @{ Row = 14; Col = 3; Value = 1.015653322938809 },
    @{ Row = 14; Col = 4; Value = 1.02314447357402 },
    @{ Row = 14; Col = 5; Value = 1.042933507076377 },
    @{ Row = 14; Col = 6; Value = 1.046733864114374 },
    @{ Row = 14; Col = 9; Value = 1.030091670737524 },
    @{ Row = 14; Col = 10; Value = 1.022371040658195 },
    @{ Row = 14; Col = 11; Value = 1.026763672824312 },
    @{ Row = 14; Col = 12; Value = 1.046478788330728 },
    @{ Row = 14; Col = 13; Value = 1.050265303431231 },
    @{ Row = 14; Col = 14; Value = 1.01156512105588 },
    @{ Row = 15; Col = 2; Value = 1.02 },
    @{ Row = 15; Col = 3; Value = 1.015743252049947 },
    @{ Row = 15; Col = 4; Value = 1.023211679183128 },
    @{ Row = 15; Col = 5; Value = 1.043045785977454 },
    @{ Row = 15; Col = 6; Value = 1.046848398045377 },
    @{ Row = 15; Col = 9; Value = 1.03010999687953 },
    @{ Row = 15; Col = 10; Value = 1.022422551527934 },
    @{ Row = 15; Col = 11; Value = 1.02681084276226 },
    @{ Row = 15; Col = 12; Value = 1.046571308474915 },
    @{ Row = 15; Col = 13; Value = 1.050360154140033 },
    @{ Row = 15; Col = 14; Value = 1.011582130137 },
    @{ Row = 16; Col = 2; Value = 1.02 },
    @{ Row = 16; Col = 3; Value = 1.016266737009357 },
    @{ Row = 16; Col = 4; Value = 1.023602742795968 },
    @{ Row = 16; Col = 5; Value = 1.043699644355094 },
    @{ Row = 16; Col = 6; Value = 1.047515266985707 },
    @{ Row = 16; Col = 9; Value = 1.030215989947695 },
    @{ Row = 16; Col = 10; Value = 1.02272218972539 },
    @{ Row = 16; Col = 11; Value = 1.027085029897313 },
    @{ Row = 16; Col = 12; Value = 1.047109938770595 },
    @{ Row = 16; Col = 13; Value = 1.050912236520142 },
    @{ Row = 16; Col = 14; Value = 1.011681067552092 },
    @{ Row = 17; Col = 2; Value = 1.02 },
    @{ Row = 17; Col = 3; Value = 1.016595156052327 },
    @{ Row = 17; Col = 4; Value = 1.023847955526079 },
    @{ Row = 17; Col = 5; Value = 1.04411009729716 },
    @{ Row = 17; Col = 6; Value = 1.047933780038103 },
    @{ Row = 17; Col = 9; Value = 1.030281883293688 },
    @{ Row = 17; Col = 10; Value = 1.022909987918611 },
    @{ Row = 17; Col = 11; Value = 1.027256700270976 },
    @{ Row = 17; Col = 12; Value = 1.047447915295739 },
    @{ Row = 17; Col = 13; Value = 1.051258551724574 },
    @{ Row = 17; Col = 14; Value = 1.011743072810253 },
    @{ Row = 18; Col = 2; Value = 1.02 },
    @{ Row = 18; Col = 3; Value = 1.01678673525541 },
    @{ Row = 18; Col = 4; Value = 1.023990950626903 },
    @{ Row = 18; Col = 5; Value = 1.044349617239778 },
    @{ Row = 18; Col = 6; Value = 1.048177964857358 },
    @{ Row = 18; Col = 9; Value = 1.030320103727611 },
    @{ Row = 18; Col = 10; Value = 1.023019470574595 },
    @{ Row = 18; Col = 11; Value = 1.027356717032351 },
    @{ Row = 18; Col = 12; Value = 1.047645090028304 },
    @{ Row = 18; Col = 13; Value = 1.051460554384491 },
    @{ Row = 18; Col = 14; Value = 1.011779219296384 },
    @{ Row = 19; Col = 2; Value = 1.02 },
    @{ Row = 19; Col = 3; Value = 1.01685206182017 },
    @{ Row = 19; Col = 4; Value = 1.024039702560184 },
    @{ Row = 19; Col = 5; Value = 1.044431305938055 },
    @{ Row = 19; Col = 6; Value = 1.048261237973251 },
    @{ Row = 19; Col = 9; Value = 1.030333099582011 },
    @{ Row = 19; Col = 10; Value = 1.023056791659635 },
    @{ Row = 19; Col = 11; Value = 1.027390800511678 },
    @{ Row = 19; Col = 12; Value = 1.047712328072383 },
    @{ Row = 19; Col = 13; Value = 1.051529432536828 },
    @{ Row = 19; Col = 14; Value = 1.011791540889294 },
    @{ Row = 20; Col = 2; Value = 1.02 },
    @{ Row = 20; Col = 3; Value = 1.016559917934227 },
    @{ Row = 20; Col = 4; Value = 1.023821649957878 },
    @{ Row = 20; Col = 5; Value = 1.044066048231519 },
    @{ Row = 20; Col = 6; Value = 1.047888869968473 },
    @{ Row = 20; Col = 9; Value = 1.030274835695434 },
    @{ Row = 20; Col = 10; Value = 1.022889844823637 },
    @{ Row = 20; Col = 11; Value = 1.027238293620247 },
    @{ Row = 20; Col = 12; Value = 1.047411649601904 },
    @{ Row = 20; Col = 13; Value = 1.05122139508003 },
    @{ Row = 20; Col = 14; Value = 1.011736422312814 },
    @{ Row = 21; Col = 2; Value = 1.02 },
    @{ Row = 21; Col = 3; Value = 1.015610343200014 },
    @{ Row = 21; Col = 4; Value = 1.023112351489183 },
    @{ Row = 21; Col = 5; Value = 1.042879850540726 },
    @{ Row = 21; Col = 6; Value = 1.046679127783769 },
    @{ Row = 21; Col = 9; Value = 1.030082900051488 },
    @{ Row = 21; Col = 10; Value = 1.022346418393409 },
    @{ Row = 21; Col = 11; Value = 1.02674112200358 },
    @{ Row = 21; Col = 12; Value = 1.046434571361995 },
    @{ Row = 21; Col = 13; Value = 1.050219970600138 },
    @{ Row = 21; Col = 14; Value = 1.011556990617815 },
    @{ Row = 22; Col = 2; Value = 1.02 },
    @{ Row = 22; Col = 3; Value = 1.015013477143862 },
    @{ Row = 22; Col = 4; Value = 1.022666097529316 },
    @{ Row = 22; Col = 5; Value = 1.042135032990148 },
    @{ Row = 22; Col = 6; Value = 1.045919179025954 },
    @{ Row = 22; Col = 9; Value = 1.029960298940605 },
    @{ Row = 22; Col = 10; Value = 1.022004238790188 },
    @{ Row = 22; Col = 11; Value = 1.02642749594322 },
    @{ Row = 22; Col = 12; Value = 1.045820596126081 },
    @{ Row = 22; Col = 13; Value = 1.049590365499986 },
    @{ Row = 22; Col = 14; Value = 1.011443995684468 },
    @{ Row = 23; Col = 2; Value = 1.02 },
    @{ Row = 23; Col = 3; Value = 1.015329870714648 },
    @{ Row = 23; Col = 4; Value = 1.022902692166664 },
    @{ Row = 23; Col = 5; Value = 1.042529780646611 },
    @{ Row = 23; Col = 6; Value = 1.04632197907291 },
    @{ Row = 23; Col = 9; Value = 1.030025474211536 },
    @{ Row = 23; Col = 10; Value = 1.02218568231273 },
    @{ Row = 23; Col = 11; Value = 1.026593852752399 },
    @{ Row = 23; Col = 12; Value = 1.046146042415228 },
    @{ Row = 23; Col = 13; Value = 1.049924127996111 },
    @{ Row = 23; Col = 14; Value = 1.011503913304248 },
    @{ Row = 24; Col = 2; Value = 1.02 },
    @{ Row = 24; Col = 3; Value = 1.016575840467665 },
    @{ Row = 24; Col = 4; Value = 1.023833536414996 },
    @{ Row = 24; Col = 5; Value = 1.044085951769272 },
    @{ Row = 24; Col = 6; Value = 1.047909162669488 },
    @{ Row = 24; Col = 9; Value = 1.030278020863123 },
    @{ Row = 24; Col = 10; Value = 1.022898946796949 },
    @{ Row = 24; Col = 11; Value = 1.027246611150881 },
    @{ Row = 24; Col = 12; Value = 1.047428036388172 },
    @{ Row = 24; Col = 13; Value = 1.051238184559453 },
    @{ Row = 24; Col = 14; Value = 1.011739427448572 },
    @{ Row = 25; Col = 2; Value = 1.02 },
    @{ Row = 25; Col = 3; Value = 1.018022418309706 },
    @{ Row = 25; Col = 4; Value = 1.024912416048646 },
    @{ Row = 25; Col = 5; Value = 1.045896102405686 },
    @{ Row = 25; Col = 6; Value = 1.049753871431435 },
    @{ Row = 25; Col = 9; Value = 1.030562688279954 },
    @{ Row = 25; Col = 10; Value = 1.023724416901122 },
    @{ Row = 25; Col = 11; Value = 1.027999559686701 },
    @{ Row = 25; Col = 12; Value = 1.048917236922605 },
    @{ Row = 25; Col = 13; Value = 1.052763183856533 },
    @{ Row = 25; Col = 14; Value = 1.012011937556457 }
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.Row, $chg.Col).Value = $chg.Value
}

